$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.422.32'
$ws.Range("E2").Value = '  -7.66%  '
$ws.Range("D3").Value = '2.902.29'
$ws.Range("E3").Value = '  -10.27%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '''473.70'
$ws.Range("E5").Value = '  -12.22%  '
$ws.Range("D6").Value = '''126.29'
$ws.Range("E6").Value = '  -7.58%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '2.903.19'
$ws.Range("E8").Value = '  -10.19%  '
$ws.Range("D9").Value = '''0.403'
$ws.Range("E9").Value = '  -12.28%  '
$ws.Range("D10").Value = '''6.68'
$ws.Range("E10").Value = '  -12.27%  '
$ws.Range("D11").Value = '''0.0969'
$ws.Range("E11").Value = '  -15.65%  '
$ws.Range("D12").Value = '''0.331'
$ws.Range("E12").Value = '  -16.36%  '
$ws.Range("E13").Value = '  -3.62%  '
$ws.Range("D14").Value = '3.392.51'
$ws.Range("E14").Value = '  -10.43%  '
$ws.Range("D15").Value = '''23.01'
$ws.Range("E15").Value = '  -11.83%  '
$ws.Range("D16").Value = '54.737.44'
$ws.Range("E16").Value = '  -7.18%  '
$ws.Range("D17").Value = '2.901.58'
$ws.Range("E17").Value = '  -10.20%  '
$ws.Range("D18").Value = '''0.0000134'
$ws.Range("E18").Value = '  -15.64%  '
$ws.Range("D19").Value = '''5.14'
$ws.Range("E19").Value = '  -13.13%  '
$ws.Range("D20").Value = '''11.52'
$ws.Range("E20").Value = '  -13.44%  '
$ws.Range("D21").Value = '''7.14'
$ws.Range("E21").Value = '  -13.99%  '
$ws.Range("D22").Value = '''308.73'
$ws.Range("E22").Value = '  -14.71%  '
$ws.Range("D24").Value = '''0.448'
$ws.Range("E24").Value = '  -14.06%  '
$ws.Range("D25").Value = '''59.13'
$ws.Range("E25").Value = '  -16.19%  '
$ws.Range("D26").Value = '''0.993'
$ws.Range("E26").Value = '  -0.74%  '
$ws.Range("D27").Value = '''0.155'
$ws.Range("E27").Value = '  -9.18%  '
$ws.Range("D28").Value = '''0.999'
$ws.Range("E28").Value = '  -0.03%  '
$ws.Range("D29").Value = '0.0₃0831'
$ws.Range("E29").Value = '  -14.99%  '
$ws.Range("D30").Value = '''6.16'
$ws.Range("E30").Value = '  -13.22%  '
$ws.Range("D31").Value = '''1.16'
$ws.Range("E31").Value = '  -5.86%  '
$ws.Range("D32").Value = '''6.16'
$ws.Range("E32").Value = '  -13.27%  '
$ws.Range("D33").Value = '''19.20'
$ws.Range("E33").Value = '  -12.71%  '
$ws.Range("D34").Value = '''1.61'
$ws.Range("E34").Value = '  -17.01%  '
$ws.Range("D35").Value = '''141.25'
$ws.Range("E35").Value = '  -13.37%  '
$ws.Range("D36").Value = '''4.25'
$ws.Range("E36").Value = '  -14.24%  '
$ws.Range("D37").Value = '''5.46'
$ws.Range("E37").Value = '  -15.12%  '
$ws.Range("D38").Value = '''1.22'
$ws.Range("E38").Value = '  -15.49%  '
$ws.Range("B39").Value = 'EnergySwap'
$ws.Range("C39").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D39").Value = '''22.50'
$ws.Range("E39").Value = '  -15.41%  '
$ws.Range("D40").Value = '''0.0621'
$ws.Range("E40").Value = '  -12.79%  '
$ws.Range("B41").Value = 'RenzoRestakedETH'
$ws.Range("C41").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D41").Value = '2.925.93'
$ws.Range("E41").Value = '  -10.39%  '
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").Value = '''1.00'
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("D43").Value = '''34.64'
$ws.Range("E43").Value = '  -15.87%  '
$ws.Range("D44").Value = '''0.967'
$ws.Range("E44").Value = '  -12.30%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").Value = '''0.599'
$ws.Range("E45").Value = '  -16.51%  '
$ws.Range("B46").Value = 'Filecoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D46").Value = '''3.44'
$ws.Range("E46").Value = '  -14.94%  '
$ws.Range("D47").Value = '''1.31'
$ws.Range("E47").Value = '  -13.37%  '
$ws.Range("D48").Value = '2.065.13'
$ws.Range("E48").Value = '  -10.20%  '
$ws.Range("D49").Value = '''18.16'
$ws.Range("E49").Value = '  -13.30%  '
$ws.Range("D50").Value = '''5.31'
$ws.Range("E50").Value = '  -15.66%  '
$ws.Range("D51").Value = '''0.0212'
$ws.Range("E51").Value = '  -12.44%  '
